$wb = $excel.ActiveWorkbook

# --- Rename headers on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add new "PO Forecast" sheet at the end (after "Monthly Trend") ---
$wsForecast = $wb.Worksheets.Add($null, $wsMonthly)
$wsForecast.Name = "PO Forecast"

# Match the page margins used on the other sheets (values are in points)
$wsForecast.PageSetup.LeftMargin = 0.75 * 72
$wsForecast.PageSetup.RightMargin = 0.75 * 72
$wsForecast.PageSetup.TopMargin = 1 * 72
$wsForecast.PageSetup.BottomMargin = 1 * 72
$wsForecast.PageSetup.HeaderMargin = 0.5 * 72
$wsForecast.PageSetup.FooterMargin = 0.5 * 72

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Copy the header style (bold, border, centered) from the existing header cells
$wsWeekly.Range("B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$data = @(
    @(45466.99999999999, 90, 84.46596963712595, 94.83906714996006),
    @(45627.99999999999, 13, 7.519291895233915, 18.12885778329256),
    @(45634.99999999999, 9, 4.008495077303803, 14.49536638426412),
    @(45641.99999999999, 6, 0.8628581221992245, 11.35841746088286),
    @(45648.99999999999, 3, -3.027310682369742, 7.536752663424573),
    @(45655.99999999999, 0, -5.686562827491834, 4.852442082196607),
    @(45662.99999999999, 0, -9.504623570949178, 1.428503477436356),
    @(45669.99999999999, 0, -12.95267576163326, -2.214279069000619),
    @(45676.99999999999, 0, -16.01494634305708, -5.679953624611109),
    @(45683.99999999999, 0, -19.47236591089338, -8.786023360335619),
    @(45690.99999999999, 0, -22.67079095273244, -12.58759680543067),
    @(45697.99999999999, 0, -26.20841525424989, -15.63828404566195)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# Copy the date-column style (format "YYYY-MM-DD HH:MM:SS") down column A
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A13").PasteSpecial(-4122)
